$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 43, shifting the old rows
# 43-46 down to 44-47 (carrying their formatting, including the date
# style on column D).
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with this week's data point.
$ws.Range("A43").Value = 8
$ws.Range("B43").Value = "Terminal La Palmera de La Serena"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44746
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100114007
$ws.Range("G43").Value = "Jengibre"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 480
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 16000
$ws.Range("M43").Value = 15500
$ws.Range("N43").Value = "$/caja 13 kilos"
$ws.Range("O43").Value = "Perú"
$ws.Range("P43").Value = 1192
$ws.Range("Q43").Value = 13
$ws.Range("R43").Value = "Hortaliza"
